# Append a new scraped listing as row 10 on the "ランサーズ" sheet, shifting the
# existing rows 10-15 down to 11-16, refresh the capture timestamp for every
# data row, widen columns B/D slightly, and rebuild the F-column hyperlinks so
# they keep tracking the right row after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- widen columns B (title) and D (price) -------------------------------
# ColumnWidth is expressed in "characters"; the engine stores the XML width
# using a +5/6 pixel-rounding offset, so back that out to land on the exact
# target values (55 / 32) the workbook should end up with.
$ws.Columns.Item(2).ColumnWidth = 55 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 32 - (5/6)

# --- insert the new row at position 10 ------------------------------------
$ws.Rows.Item(10).Insert()

# --- populate the newly inserted row 10 ------------------------------------
$newTimestamp = "2025-11-28 01:48:43"

$ws.Range("A10").Value = $newTimestamp
$ws.Range("B10").Value = "初回 2026年1月創業 コンサル会社のバックオフィス業務フロー設計・マニュアル化、IT導入 一括見積依頼"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5442904"
$ws.Range("G10").Value = 55
$ws.Range("H10").Value = "◆コンサル"

# --- refresh the capture timestamp on every other data row -----------------
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A8").Value = $newTimestamp
$ws.Range("A9").Value = $newTimestamp
$ws.Range("A11").Value = $newTimestamp
$ws.Range("A12").Value = $newTimestamp
$ws.Range("A13").Value = $newTimestamp
$ws.Range("A14").Value = $newTimestamp
$ws.Range("A15").Value = $newTimestamp
$ws.Range("A16").Value = $newTimestamp

# --- rebuild the F-column hyperlinks ---------------------------------------
# Row-insert does not renumber the sheet's <hyperlinks> table in this engine,
# so wipe the lot and re-add them in row order; that mints rId1..rId15
# pointing at the right targets for the new row layout. Re-apply the
# "Hyperlink" cell style afterwards since Hyperlinks.Add mints its own
# (functionally identical) style record otherwise.
$ws.Range("F2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5427956")
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5442448")
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5441932")
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5442360")
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5442482")
$ws.Range("F6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5442416")
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5441612")
$ws.Range("F8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5442063")
$ws.Range("F9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5442904")
$ws.Range("F10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5442064")
$ws.Range("F11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5442625")
$ws.Range("F12").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5442153")
$ws.Range("F13").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5442106")
$ws.Range("F14").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5442432")
$ws.Range("F15").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5442169")
$ws.Range("F16").Style = "Hyperlink"
